# Add pending-approval rows 104-107 (new vendor entries appended to PENDING_APPROVAL_SHEET)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104
$ws.Range('A104').Value = 'WGG 02'
$ws.Range('B104').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('C104').Value = 46328
$ws.Range('C104').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('D104').Value = 286962
$ws.Range('E104').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('F104').Value = 34413429360
$ws.Range('G104').Value = 'NEFT'
$ws.Range('H104').Value = 'SBIN0003229'
$ws.Range('I104').Value = 'AAAFW8862C'
$ws.Range('J104').Value = '32AAAFW8862C1Z9'
$ws.Range('L104').Value = '9745dde2-4798-4594-8aa3-f4638d15e9f7'
$ws.Range('U104').Value = 'pending'
$ws.Range('V104').Value = 0
$ws.Range('X104').Value = '0 RPA_ID : 28134e0fd0'
$ws.Range('Y104').Value = 0
$ws.Range('Z104').Value = 0
$ws.Range('AA104').Value = 'officeadmin@westernidc.com'
$ws.Range('AB104').Value = 'ESTIMATION NOT MATCHED'
$ws.Range('AC104').Value = 0
$ws.Range('AD104').Value = 0
$ws.Range('AE104').Value = 0

# Row 105
$ws.Range('A105').Value = 'WGE 303'
$ws.Range('B105').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('C105').Value = 46328
$ws.Range('C105').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('D105').Value = 286962
$ws.Range('E105').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('F105').Value = 34413429360
$ws.Range('G105').Value = 'NEFT'
$ws.Range('H105').Value = 'SBIN0003229'
$ws.Range('I105').Value = 'AAAFW8862C'
$ws.Range('J105').Value = '32AAAFW8862C1Z9'
$ws.Range('K105').Value = 'Hemilda Benson'
$ws.Range('L105').Value = '764da05a-eef7-4e4b-a23f-0789be25f190'
$ws.Range('U105').Value = 'pending'
$ws.Range('V105').Value = 27796.23
$ws.Range('X105').Value = 'Closing ICICI credit card RPA_ID : 0d5cf2c164'
$ws.Range('Y105').Value = 0
$ws.Range('Z105').Value = 'PAYMENT'
$ws.Range('AA105').Value = 'officeadmin@westernidc.com'
$ws.Range('AB105').Value = 'ESTIMATION NOT MATCHED'
$ws.Range('AC105').Value = 0
$ws.Range('AD105').Value = 0
$ws.Range('AE105').Value = 0

# Row 106
$ws.Range('A106').Value = 'WGG 02'
$ws.Range('B106').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('C106').Value = 46328
$ws.Range('C106').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('D106').Value = 286962
$ws.Range('E106').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('F106').Value = 34413429360
$ws.Range('G106').Value = 'NEFT'
$ws.Range('H106').Value = 'SBIN0003229'
$ws.Range('I106').Value = 'AAAFW8862C'
$ws.Range('J106').Value = '32AAAFW8862C1Z9'
$ws.Range('L106').Value = '16703548-9f1a-4d4d-a116-49db98e90a68'
$ws.Range('U106').Value = 'pending'
$ws.Range('V106').Value = 0
$ws.Range('X106').Value = '0 RPA_ID : 45413f950a'
$ws.Range('Y106').Value = 0
$ws.Range('Z106').Value = 0
$ws.Range('AA106').Value = 'officeadmin@westernidc.com'
$ws.Range('AB106').Value = 'ESTIMATION NOT MATCHED'
$ws.Range('AC106').Value = 0
$ws.Range('AD106').Value = 0
$ws.Range('AE106').Value = 0

# Row 107
$ws.Range('A107').Value = 'WGA009'
$ws.Range('B107').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('C107').Value = 46328
$ws.Range('C107').NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range('D107').Value = 286962
$ws.Range('E107').Value = 'Western Interior Designers & Marine Contractors'
$ws.Range('F107').Value = 34413429360
$ws.Range('G107').Value = 'NEFT'
$ws.Range('H107').Value = 'SBIN0003229'
$ws.Range('I107').Value = 'AAAFW8862C'
$ws.Range('J107').Value = '32AAAFW8862C1Z9'
$ws.Range('K107').Value = 'SNAPDRY MAX LLP'
$ws.Range('L107').Value = '0a7faa6e-79e2-4f79-adac-728b0d527f46'
$ws.Range('M107').Value = 13160200032800
$ws.Range('N107').Value = 'FDRL0001316'
$ws.Range('U107').Value = 'pending'
$ws.Range('V107').Value = 3354
$ws.Range('X107').Value = 'wash and fold carpets on Jan RPA_ID : 21364abc02'
$ws.Range('Y107').Value = 0
$ws.Range('Z107').Value = 'PAYMENT'
$ws.Range('AA107').Value = 'officeadmin@westernidc.com'
$ws.Range('AB107').Value = 'ESTIMATION NOT MATCHED'
$ws.Range('AC107').Value = 0
$ws.Range('AD107').Value = 0
$ws.Range('AE107').Value = 0
